$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp text (shared string used in A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 16:18"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4815895
$ws.Range("C4").Value = 2248
$ws.Range("D4").Value = 2380584
$ws.Range("E4").Value = 2276935
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 158376

# Row 6 - India
$ws.Range("B6").Value = 1830949
$ws.Range("C6").Value = 26247
$ws.Range("D6").Value = 1200303
$ws.Range("E6").Value = 592161
$ws.Range("G6").Value = 324
$ws.Range("H6").Value = 38485

# Row 22 - Argentina
$ws.Range("D22").Value = 91302
$ws.Range("E22").Value = 106950
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 3667

# Row 31 - Ecuador
$ws.Range("D31").Value = 59344
$ws.Range("E31").Value = 21430

# Row 33 - Suecia
$ws.Range("B33").Value = 81012
$ws.Range("C33").Value = 43
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 5744

# Row 59 - Azerbaiyan
$ws.Range("B59").Value = 32684
$ws.Range("C59").Value = 241
$ws.Range("D59").Value = 27760
$ws.Range("E59").Value = 4456
$ws.Range("G59").Value = 6
$ws.Range("H59").Value = 468

# Row 61 - Serbia
$ws.Range("B61").Value = 26451
$ws.Range("C61").Value = 258
$ws.Range("E61").Value = 11806
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 598

# Row 80 - Bosnia y Herzegovina
$ws.Range("B80").Value = 12462
$ws.Range("C80").Value = 166
$ws.Range("D80").Value = 6359
$ws.Range("E80").Value = 5741
$ws.Range("G80").Value = 10
$ws.Range("H80").Value = 362

# Row 86 - Noruega
$ws.Range("B86").Value = 9312
$ws.Range("C86").Value = 44
$ws.Range("E86").Value = 304

# Row 189 - Papua Nueva Guinea
$ws.Range("B189").Value = 111
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 75
